$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.749.68"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.38%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.634.19"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.80%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "537.29"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.63%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.13"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +3.58%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.32%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.56"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.16%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.28%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.66%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.64%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.101.22"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.99%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "59.643.82"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.14%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.93"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.22%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.625.91"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.93%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000135"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.56%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "342.11"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.71%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.73%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.43%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.38"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.50%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.50"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.38%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.84%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.25%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.68%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +5.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +4.07%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.98"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.59%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "150.90"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.34%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.65%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.17%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.72%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.836"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.76%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.823"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.92%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "288.99"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +8.09%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.86%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.998"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.08%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.73"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.20%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.24%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0531"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +3.67%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.967.20"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.51%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.62%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.44"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.78%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.61%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "111.07"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.49%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.73"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.28%  "
